# Remove the "P_kbar" column (column S) from the worksheet.
# Excel shifts the remaining columns (Teq15_2H2O, Teq15_4H2O) left by one,
# which matches the diff: P_kbar is dropped from sharedStrings and the
# old column S data is removed while columns T/U slide into S/T.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("S1").EntireColumn.Delete()
